# Updates the cryptos list sheet with refreshed price/volume data
# (row 41 <-> 42 coin swap: Filecoin now ranks above EnergySwap)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceCell {
    param($Cell, $Text, $ForceText)
    if ($ForceText) {
        $Cell.NumberFormat = "@"
    }
    $Cell.Value = $Text
}

# Rows whose Price (D) and/or Volume(1h) (E) values changed.
# ForceText marks values that look like plain numbers (single decimal point)
# so Excel doesn't silently coerce them to floating point numbers and lose
# formatting (trailing zeros, exact digits, etc.).
$updates = @(
    @{ Row=2; D='67.499.80'; E='  +3.72%  '; ForceText=$false },
    @{ Row=3; D='3.257.42'; E='  +2.98%  '; ForceText=$false },
    @{ Row=4; D=$null; E='  -0.07%  '; ForceText=$false },
    @{ Row=5; D='578.52'; E='  +1.69%  '; ForceText=$true },
    @{ Row=6; D='182.30'; E='  +7.16%  '; ForceText=$true },
    @{ Row=7; D=$null; E='  -0.08%  '; ForceText=$false },
    @{ Row=8; D='0.597'; E='  -3.10%  '; ForceText=$true },
    @{ Row=9; D='3.253.31'; E='  +3.00%  '; ForceText=$false },
    @{ Row=10; D=$null; E='  +6.13%  '; ForceText=$false },
    @{ Row=11; D='6.75'; E='  +3.10%  '; ForceText=$true },
    @{ Row=12; D=$null; E='  +5.52%  '; ForceText=$false },
    @{ Row=13; D='3.815.39'; E='  +2.64%  '; ForceText=$false },
    @{ Row=14; D=$null; E='  +1.15%  '; ForceText=$false },
    @{ Row=15; D='28.49'; E='  +5.15%  '; ForceText=$true },
    @{ Row=16; D='67.467.41'; E='  +3.79%  '; ForceText=$false },
    @{ Row=17; D='0.0000168'; E='  +3.56%  '; ForceText=$true },
    @{ Row=18; D='3.254.07'; E='  +2.80%  '; ForceText=$false },
    @{ Row=19; D='5.85'; E='  +2.12%  '; ForceText=$true },
    @{ Row=20; D='13.53'; E='  +5.67%  '; ForceText=$true },
    @{ Row=21; D='376.28'; E='  +5.26%  '; ForceText=$true },
    @{ Row=22; D='7.63'; E='  +4.74%  '; ForceText=$true },
    @{ Row=23; D=$null; E='  -0.22%  '; ForceText=$false },
    @{ Row=24; D='71.24'; E='  +2.95%  '; ForceText=$true },
    @{ Row=25; D='0.512'; E='  +3.04%  '; ForceText=$true },
    @{ Row=26; D=$null; E='  +3.05%  '; ForceText=$false },
    @{ Row=27; D='9.58'; E='  -2.04%  '; ForceText=$true },
    @{ Row=28; D='0.180'; E='  +2.05%  '; ForceText=$true },
    @{ Row=29; D=$null; E='  +0.37%  '; ForceText=$false },
    @{ Row=30; D=$null; E='  +8.55%  '; ForceText=$false },
    @{ Row=31; D=$null; E='  +3.21%  '; ForceText=$false },
    @{ Row=32; D='22.71'; E='  +3.25%  '; ForceText=$true },
    @{ Row=33; D=$null; E='  +0.01%  '; ForceText=$false },
    @{ Row=34; D=$null; E='  +5.59%  '; ForceText=$false },
    @{ Row=35; D='6.93'; E='  +4.48%  '; ForceText=$true },
    @{ Row=36; D='163.75'; E='  +4.16%  '; ForceText=$true },
    @{ Row=37; D=$null; E='  +4.02%  '; ForceText=$false },
    @{ Row=38; D='0.849'; E='  +1.64%  '; ForceText=$true },
    @{ Row=39; D=$null; E='  +4.60%  '; ForceText=$false },
    @{ Row=40; D=$null; E='  +12.97%  '; ForceText=$false },
    @{ Row=43; D='2.62'; E='  +5.90%  '; ForceText=$true },
    @{ Row=44; D='357.93'; E='  +10.68%  '; ForceText=$true },
    @{ Row=45; D='2.713.30'; E='  +1.73%  '; ForceText=$false },
    @{ Row=46; D='25.46'; E='  +5.84%  '; ForceText=$true },
    @{ Row=47; D='40.84'; E='  +3.30%  '; ForceText=$true },
    @{ Row=48; D=$null; E='  +3.96%  '; ForceText=$false },
    @{ Row=49; D='0.0280'; E='  +2.78%  '; ForceText=$true },
    @{ Row=50; D='1.00'; E='  +6.25%  '; ForceText=$true },
    @{ Row=51; D=$null; E='  -1.07%  '; ForceText=$false }
)

foreach ($item in $updates) {
    if ($item.D -ne $null) {
        Set-PriceCell ($ws.Cells.Item($item.Row, 4)) $item.D $item.ForceText
    }
    Set-PriceCell ($ws.Cells.Item($item.Row, 5)) $item.E $false
}

# Row 41 and 42 swapped places: Filecoin moved up (now rank 39), EnergySwap
# moved down (now rank 40). Update Coin, Link, Price and Volume(1h) columns
# for both rows to reflect the new order and values.
$ws.Cells.Item(41, 2).Value = "Filecoin"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-PriceCell ($ws.Cells.Item(41, 4)) "4.65" $true
Set-PriceCell ($ws.Cells.Item(41, 5)) "  +11.23%  " $false

$ws.Cells.Item(42, 2).Value = "EnergySwap"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-PriceCell ($ws.Cells.Item(42, 4)) "26.63" $true
Set-PriceCell ($ws.Cells.Item(42, 5)) "  +1.98%  " $false
